{"js": "// Apply the LOQ4055 course-page text updates described by the diff.\n// Each entry is [oldText, newText]; we search the body for an exact\n// (case-sensitive) match and replace it in place.\nconst replacements = [\n  [\n    \"Ativa\u00e7\u00e3o: 01/01/2022\",\n    \"Ativa\u00e7\u00e3o: 01/01/2024\"\n  ],\n  [\n    \"Curso (semestre ideal): EQD (3), EQN (4)\",\n    \"Curso (semestre ideal): EQN (4)\"\n  ],\n  [\n    \"Fornecer aos alunos conceitos fundamentais para compreens\u00e3o da Qu\u00edmica Inorg\u00e2nica por meio da experimenta\u00e7\u00e3o, desenvolvendo a capacidade de realizarem pr\u00e1ticas no laborat\u00f3rio que estimulem o seu pensamento cient\u00edfico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de car\u00e1ter inorg\u00e2nico com interesse industrial.\",\n    \"Fornecer aos alunos conceitos fundamentos para a compreens\u00e3o da Qu\u00edmica Inorg\u00e2nica, de forma a capacit\u00e1-lo a descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de car\u00e1ter inorg\u00e2nico com interesse industrial.\"\n  ],\n  [\n    \"Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest.\",\n    \"Provide students with fundamental concepts for understanding Inorganic Chemistry, in order to enable them to describe and interpret the properties of elements and their compounds, especially those of an inorganic nature with industrial interest\"\n  ],\n  [\n    \"Compostos de Coordena\u00e7\u00e3o. Materiais inorg\u00e2nicos de interesse industrial. Purifica\u00e7\u00e3o e Identifica\u00e7\u00e3o de Compostos Inorg\u00e2nicos. S\u00edntese de sais e obten\u00e7\u00e3o de Compostos de Alum\u00ednio.\",\n    \"M\u00e9todos de separa\u00e7\u00e3o e obten\u00e7\u00e3o dos elementos, extra\u00e7\u00e3o mineral. Hidrog\u00eanio. Metais alcalinos. Metais alcalino terrosos. Alum\u00ednio. Metais de transi\u00e7\u00e3o. Compostos de coordena\u00e7\u00e3o. Halog\u00eanios.\"\n  ],\n  [\n    \"Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.\",\n    \"Methods for separating and obtaining elements, mineral extraction. Hydrogen. Alkaline metals. Alkaline earth metals. Aluminum. Transition metals. Coordination compounds. Halogens.\"\n  ],\n  [\n    \"Compostos de Coordena\u00e7\u00e3o: Estrutura, liga\u00e7\u00f5es, rea\u00e7\u00f5es e aplica\u00e7\u00f5es. Exemplos e aplica\u00e7\u00f5es de materiais inorg\u00e2nicos de interesse industrial. S\u00ednteses: Sal Simples, Sal Duplo e Sal Complexo. Prepara\u00e7\u00e3o de Compostos de Alum\u00ednio.\",\n    \"M\u00e9todos de separa\u00e7\u00e3o e obten\u00e7\u00e3o dos elementos, extra\u00e7\u00e3o mineral. Propriedades, obten\u00e7\u00e3o e aplica\u00e7\u00f5es dos seguintes elementos/grupos e seus compostos: Hidrog\u00eanio; Metais alcalinos (ind\u00fastria cloro-\u00e1lcali; processo Solvay); Metais alcalino terrosos; Alum\u00ednio (processo Bayer); Metais de transi\u00e7\u00e3o; Compostos de coordena\u00e7\u00e3o e Halog\u00eanios.\"\n  ],\n  [\n    \"Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.\",\n    \"Methods of separating and obtaining the elements, mineral extraction. Properties, obtaining and applications of the following elements/groups and their compounds: Hydrogen; Alkali metals (chlor-alkali industry; Solvay process); Alkaline earth metals; Aluminum (Bayer process); Transition metals; Coordination compounds and Halogens.\"\n  ],\n  [\n    \"Ser\u00e3o oferecidas aulas expositivas e pr\u00e1ticas.\",\n    \"Ser\u00e3o oferecidas aulas expositivas.\"\n  ],\n  [\n    \"CHANG, Raymond. Qu\u00edmica geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Qu\u00edmica a ci\u00eancia central. 9.ed. S\u00e3o Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Qu\u00edmica geral. Rio de Janeiro: Ed. Livros T\u00e9cnicos Cient\u00edficos, 1981.LEE, J. D., tradu\u00e7\u00e3o Qu\u00edmica Inorg\u00e2nica n\u00e3o t\u00e3o concisa da 5\u00aa edi\u00e7\u00e3o inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Qu\u00edmica Inorg\u00e2nica tradu\u00e7\u00e3o da 4\u00aa edi\u00e7\u00e3o. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Qu\u00edmica - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3\u00aa ed., 1973.\",\n    \"WELLER, Mark; OVERTON, Tina; ROURKE, Jonathan; et al. Qu\u00edmica inorg\u00e2nica. Porto Alegre, Bookman, 6\u00aa Ed, 2017. E-book. CHANG, Raymond. Qu\u00edmica geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Qu\u00edmica a ci\u00eancia central. 9.ed. S\u00e3o Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Qu\u00edmica geral. Rio de Janeiro: Ed. Livros T\u00e9cnicos Cient\u00edficos, 1981.LEE, J. D., tradu\u00e7\u00e3o Qu\u00edmica Inorg\u00e2nica n\u00e3o t\u00e3o concisa da 5\u00aa edi\u00e7\u00e3o inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Qu\u00edmica Inorg\u00e2nica tradu\u00e7\u00e3o da 4\u00aa edi\u00e7\u00e3o. Editora Bookman, Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Qu\u00edmica - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3\u00aa ed., 1973.\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText.substring(0, 60));\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the LOQ4055 course-page text updates described by the diff.\n# Each entry is a (FindText, ReplaceText) pair; we run Find/Replace\n# across the whole document content for each one.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"Ativa\u00e7\u00e3o: 01/01/2022\", \"Ativa\u00e7\u00e3o: 01/01/2024\"),\n    @(\"Curso (semestre ideal): EQD (3), EQN (4)\", \"Curso (semestre ideal): EQN (4)\"),\n    @(\n        \"Fornecer aos alunos conceitos fundamentais para compreens\u00e3o da Qu\u00edmica Inorg\u00e2nica por meio da experimenta\u00e7\u00e3o, desenvolvendo a capacidade de realizarem pr\u00e1ticas no laborat\u00f3rio que estimulem o seu pensamento cient\u00edfico. Descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de car\u00e1ter inorg\u00e2nico com interesse industrial.\",\n        \"Fornecer aos alunos conceitos fundamentos para a compreens\u00e3o da Qu\u00edmica Inorg\u00e2nica, de forma a capacit\u00e1-lo a descrever e interpretar as propriedades dos elementos e de seus compostos, principalmente daqueles de car\u00e1ter inorg\u00e2nico com interesse industrial.\"\n    ),\n    @(\n        \"Provide the students with fundamental concepts for understanding Inorganic Chemistry by means of experimentation, developing the ability to carry out practices in the laboratory which stimulate their scientific thinking. Describe and interpret the properties of the elements and their compounds, especially those of an inorganic feature with industrial interest.\",\n        \"Provide students with fundamental concepts for understanding Inorganic Chemistry, in order to enable them to describe and interpret the properties of elements and their compounds, especially those of an inorganic nature with industrial interest\"\n    ),\n    @(\n        \"Compostos de Coordena\u00e7\u00e3o. Materiais inorg\u00e2nicos de interesse industrial. Purifica\u00e7\u00e3o e Identifica\u00e7\u00e3o de Compostos Inorg\u00e2nicos. S\u00edntese de sais e obten\u00e7\u00e3o de Compostos de Alum\u00ednio.\",\n        \"M\u00e9todos de separa\u00e7\u00e3o e obten\u00e7\u00e3o dos elementos, extra\u00e7\u00e3o mineral. Hidrog\u00eanio. Metais alcalinos. Metais alcalino terrosos. Alum\u00ednio. Metais de transi\u00e7\u00e3o. Compostos de coordena\u00e7\u00e3o. Halog\u00eanios.\"\n    ),\n    @(\n        \"Coordination Compounds. Inorganic materials of industrial interest. Purification and Identification of Inorganic Compounds. Synthesis: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.\",\n        \"Methods for separating and obtaining elements, mineral extraction. Hydrogen. Alkaline metals. Alkaline earth metals. Aluminum. Transition metals. Coordination compounds. Halogens.\"\n    ),\n    @(\n        \"Compostos de Coordena\u00e7\u00e3o: Estrutura, liga\u00e7\u00f5es, rea\u00e7\u00f5es e aplica\u00e7\u00f5es. Exemplos e aplica\u00e7\u00f5es de materiais inorg\u00e2nicos de interesse industrial. S\u00ednteses: Sal Simples, Sal Duplo e Sal Complexo. Prepara\u00e7\u00e3o de Compostos de Alum\u00ednio.\",\n        \"M\u00e9todos de separa\u00e7\u00e3o e obten\u00e7\u00e3o dos elementos, extra\u00e7\u00e3o mineral. Propriedades, obten\u00e7\u00e3o e aplica\u00e7\u00f5es dos seguintes elementos/grupos e seus compostos: Hidrog\u00eanio; Metais alcalinos (ind\u00fastria cloro-\u00e1lcali; processo Solvay); Metais alcalino terrosos; Alum\u00ednio (processo Bayer); Metais de transi\u00e7\u00e3o; Compostos de coordena\u00e7\u00e3o e Halog\u00eanios.\"\n    ),\n    @(\n        \"Coordination Compounds: Structure, bonds, reactions and applications. Examples and applications of industrial interest inorganic materials. Syntheses: Simple Salt, Double Salt and Complex Salt. Preparation of Aluminum Compounds.\",\n        \"Methods of separating and obtaining the elements, mineral extraction. Properties, obtaining and applications of the following elements/groups and their compounds: Hydrogen; Alkali metals (chlor-alkali industry; Solvay process); Alkaline earth metals; Aluminum (Bayer process); Transition metals; Coordination compounds and Halogens.\"\n    ),\n    @(\"Ser\u00e3o oferecidas aulas expositivas e pr\u00e1ticas.\", \"Ser\u00e3o oferecidas aulas expositivas.\"),\n    @(\n        \"CHANG, Raymond. Qu\u00edmica geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Qu\u00edmica a ci\u00eancia central. 9.ed. S\u00e3o Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Qu\u00edmica geral. Rio de Janeiro: Ed. Livros T\u00e9cnicos Cient\u00edficos, 1981.LEE, J. D., tradu\u00e7\u00e3o Qu\u00edmica Inorg\u00e2nica n\u00e3o t\u00e3o concisa da 5\u00aa edi\u00e7\u00e3o inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Qu\u00edmica Inorg\u00e2nica tradu\u00e7\u00e3o da 4\u00aa edi\u00e7\u00e3o. Editora Bookman,Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Qu\u00edmica - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3\u00aa ed., 1973.\",\n        \"WELLER, Mark; OVERTON, Tina; ROURKE, Jonathan; et al. Qu\u00edmica inorg\u00e2nica. Porto Alegre, Bookman, 6\u00aa Ed, 2017. E-book. CHANG, Raymond. Qu\u00edmica geral: conceitos essenciais. 4.ed. s.l.:Ed. AMGH Editora Ltda., 2010.BROWN, T.L. ET al. Qu\u00edmica a ci\u00eancia central. 9.ed. S\u00e3o Paulo: Pearson Prentice Hall, 2005-2007.BRADY, J ; HUMISTON, G.E. Qu\u00edmica geral. Rio de Janeiro: Ed. Livros T\u00e9cnicos Cient\u00edficos, 1981.LEE, J. D., tradu\u00e7\u00e3o Qu\u00edmica Inorg\u00e2nica n\u00e3o t\u00e3o concisa da 5\u00aa edi\u00e7\u00e3o inglesa. Editora Edgard Blucher Ltda. SP-2001.SHRIVER, D. e ATKINS, P. Qu\u00edmica Inorg\u00e2nica tradu\u00e7\u00e3o da 4\u00aa edi\u00e7\u00e3o. Editora Bookman, Porto Alegre-RS, 2008.QUAGLIANO, J.V; VALLARINO, L.M. Qu\u00edmica - Ed. Guanabara Koogan S.A. - Rio de Janeiro - 3\u00aa ed., 1973.\"\n    )\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, $wdReplaceAll)\n}\n"}
